$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "#constant_species" parameter was previously disabled (leading "#") and
# left without an input value. Re-enable it by dropping the "#" prefix.
$ws.Range("A14").Value = "constant_species"

# Match the formatting already used by the other populated cells in column C
# (10pt Arial, left/vertical-centered) instead of the placeholder style that
# was only in use while C14 was still empty.
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the constant-species list that had been missing.
$ws.Range("C14").Value = "HO2, N, CO2, H2O, CO, O2, N2, OH, O, H2, H, O3, "

# Reflect that C14 is where editing/selection ended up when the file was saved.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C14").Select()
